$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.952294
$ws.Range("N2").Value = 5.856882
$ws.Range("O2").Value = 0.07575070565202183
$ws.Range("P2").Value = 0.07575070565202184
$ws.Range("Q2").Value = 0.3078214488033333
$ws.Range("R2").Value = 2.77039303923
$ws.Range("S2").Value = 0.07575070565202183
$ws.Range("T2").Value = 0.07575070565202184

# Row 3
$ws.Range("O3").Value = 0.5679402069281436
$ws.Range("P3").Value = 0.5679402069281437
$ws.Range("S3").Value = 0.5679402069281436
$ws.Range("T3").Value = 0.5679402069281437

# Row 4
$ws.Range("M4").Value = 9.009963000000001
$ws.Range("N4").Value = 27.029889
$ws.Range("O4").Value = 0.3495944028658632
$ws.Range("P4").Value = 0.3495944028658634
$ws.Range("Q4").Value = 1.420615882815
$ws.Range("R4").Value = 12.785542945335
$ws.Range("S4").Value = 0.3495944028658632
$ws.Range("T4").Value = 0.3495944028658634

# Row 5
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.173055
$ws.Range("N5").Value = 0.519165
$ws.Range("O5").Value = 0.006714684553971194
$ws.Range("P5").Value = 0.006714684553971196
$ws.Range("Q5").Value = 0.027285870275
$ws.Range("R5").Value = 0.245572832475
$ws.Range("S5").Value = 0.006714684553971194
$ws.Range("T5").Value = 0.006714684553971196
